$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sat right after ", BEREKET HAILE"
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Revision Date: 04/16/2017 -> 05/04/2017
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Revision Date:", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $ptext = $para.Range.Text
    $pstart = $para.Range.Start
    $idx = $ptext.IndexOf("04/16/2017")
    if ($idx -ge 0) {
        $absStart = $pstart + $idx
        # "16" -> "04"  (update this one first so the earlier offset stays valid)
        $rDay = $d.Range($absStart + 3, $absStart + 5)
        $rDay.Text = "04"
        # "4" -> "5"
        $rMonth = $d.Range($absStart + 1, $absStart + 2)
        $rMonth.Text = "5"
    }
}

# ---------------------------------------------------------------------------
# 3) Rework the Product Backlog list:
#    - split the long "As a Developer, I want to collect data..." paragraph
#      into its own (slightly re-punctuated) item plus seven new list items
#    - drop the old manually numbered "6." / "7." paragraphs
#    - move the "_GoBack" bookmark onto the end of the new last list item
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$searchText = "As a Developer, I want to collect data in the Android App in a centralized, structured way in      order to allow for the back-end to have an easier time processing it."
$found2 = $rng2.Find.Execute($searchText, $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)

if ($found2) {
    $newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Developer, I want to collect data in the Android App in a centralized, structured way in      order to allow for the back-end to have an easier time processing it</w:t></w:r>
  <w:r><w:t xml:space="preserve">. </w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Product Sponsor, I want to be able to view all users of the Admin Panel categorically.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Developer, I want to be able to validate login information in the Android login system in      order to provide for a more accurate login system.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a User, I want to have a pleasing Chrome Extension UI to interact with and look at in order to have a pleasing experience when using the extension.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Product Sponsor, I want the Admin Panel presented in an aesthetically pleasing way so that everything is structured and clear when I visit a page.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Developer, I want the Android Application data to be buffered so that when data is sent to the backend it is done so in a clear, concise manner.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r><w:t>As a Developer, I want to debug and test the Android Application in order to reveal any bugs or issues that present itself.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:suppressAutoHyphens/>
    <w:spacing w:line="100" w:lineRule="atLeast"/>
    <w:ind w:hanging="360"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>As a User, I want the Android Application to have a nice user interface in order to provide a better viewing and user experience.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rng2.InsertXML($newParasXml)
}

# ---------------------------------------------------------------------------
# 4) Remove the text from the old "6." item (now an empty paragraph) and
#    delete the whole old "7." paragraph entirely.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("6.    As a Product Sponsor, I want to be able to view all users of the Admin Panel categorically.", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Text = ""
}

$rng4 = $d.Content
$found4 = $rng4.Find.Execute("7.     As a Developer, I want to be able to validate login information in the Android login system in      order to provide for a more accurate login system.", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found4) {
    $p4 = $rng4.Paragraphs(1)
    $pRange = $d.Range($p4.Range.Start, $p4.Range.End)
    $pRange.Text = ""
}

# ---------------------------------------------------------------------------
# 5) Re-add the "_GoBack" bookmark at the end of the new last list item
#    ("As a User, I want the Android Application to have a nice user
#    interface ...")
# ---------------------------------------------------------------------------
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("As a User, I want the Android Application to have a nice user interface in order to provide a better viewing and user experience.", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found5) {
    $bmRange = $d.Range($rng5.End, $rng5.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "done"
